# Update the Australia MSME summary figures to their more precise values.
#
# The target cells currently hold their numbers as literal TEXT (shared
# strings), not as numeric values. A plain `Range.Value = "22.67"` assignment
# would have Excel auto-detect the numeric-looking string and store it as a
# Number, which does not match the source data (and also drags in binary
# floating point noise like 22.670000000000002). Prefixing the literal with
# a leading apostrophe forces Excel to keep it as Text, and re-applying the
# "Normal" cell style afterwards clears the transient quote-prefix styling
# so the cell's style index is left exactly as it was before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enterprises density (per 1000 people): Micro / SMEs / MSMEs
$ws.Range("B11").Value = "'22.67"
$ws.Range("B11").Style = "Normal"

$ws.Range("C11").Value = "'13.84"
$ws.Range("C11").Style = "Normal"

$ws.Range("D11").Value = "'36.51"
$ws.Range("D11").Style = "Normal"

# Enterprises (% of total): Micro / SMEs / MSMEs
$ws.Range("B13").Value = "'61.62"
$ws.Range("B13").Style = "Normal"

$ws.Range("C13").Value = "'37.61"
$ws.Range("C13").Style = "Normal"

$ws.Range("D13").Value = "'99.23"
$ws.Range("D13").Style = "Normal"
